# Applies the "typenumeric" action edit to sheet TC2:
#  - B3 changes from "randemail" to the new value "typenumeric"
#  - H3 gets the numeric value 345434
#  - Selection moves to H3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC2")

$ws.Range("B3").Value = "typenumeric"
$ws.Range("H3").Value = 345434

$ws.Activate()
$ws.Range("H3").Select()
